# Generate Report for Handoff
# The file "c2b6a63f-fef9-4b2e-9b7c-e407c6336a56.md" has been re-handed-off:
#   - its Status changes from "Handed back: in sync with en-US" to "Ready for handoff"
#     in the Overview sheet (both language columns) and in each language sheet.
#   - the Latest Handoff Date/Datetime is updated to a new timestamp for each language.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# ---- Overview sheet ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus
$overview.Range("D3").Value = "2016-37-18 00:37:48"

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("E3").Value = "2016-03-18 00:37:44"

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("E3").Value = "2016-03-18 00:37:48"
